# Update the math-drill worksheet table cells to the new problems.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$failures = 0
$cell = $t.Cell(1, 1)
$ok = $cell.Range.Find.Execute("35÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷4=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(1, 2)
$ok = $cell.Range.Find.Execute("28÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷9=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(1, 3)
$ok = $cell.Range.Find.Execute("63÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷8=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(1, 4)
$ok = $cell.Range.Find.Execute("40÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷2=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(1, 5)
$ok = $cell.Range.Find.Execute("26÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷7=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(5, 1)
$ok = $cell.Range.Find.Execute("69÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷2=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(5, 2)
$ok = $cell.Range.Find.Execute("20÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷6=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(5, 3)
$ok = $cell.Range.Find.Execute("46÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷9=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(5, 4)
$ok = $cell.Range.Find.Execute("19÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷2=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(5, 5)
$ok = $cell.Range.Find.Execute("12÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷7=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(9, 1)
$ok = $cell.Range.Find.Execute("34÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷4=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(9, 2)
$ok = $cell.Range.Find.Execute("76÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷8=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(9, 3)
$ok = $cell.Range.Find.Execute("50÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷9=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(9, 4)
$ok = $cell.Range.Find.Execute("55÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷7=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(9, 5)
$ok = $cell.Range.Find.Execute("17÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷6=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(13, 1)
$ok = $cell.Range.Find.Execute("93÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷9=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(13, 2)
$ok = $cell.Range.Find.Execute("99÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷2=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(13, 3)
$ok = $cell.Range.Find.Execute("60÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷5=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(13, 4)
$ok = $cell.Range.Find.Execute("49÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷5=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(13, 5)
$ok = $cell.Range.Find.Execute("13÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷2=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(17, 1)
$ok = $cell.Range.Find.Execute("82÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷6=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(17, 2)
$ok = $cell.Range.Find.Execute("36÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷5=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(17, 3)
$ok = $cell.Range.Find.Execute("53÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷4=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(17, 4)
$ok = $cell.Range.Find.Execute("67÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷5=", 2)
if (-not $ok) { $failures++ }
$cell = $t.Cell(17, 5)
$ok = $cell.Range.Find.Execute("32÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷8=", 2)
if (-not $ok) { $failures++ }
Write-Output ("Replacements failed: " + $failures)
